$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89 (pushes existing rows 89..209 down to 90..210)
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new record
$ws.Range("A89").Value = 5
$ws.Range("B89").Value = "Macroferia Regional de Talca"
$ws.Range("C89").Value = "Maule"
$ws.Range("D89").Value = 44771
$ws.Range("E89").Value = 7
$ws.Range("F89").Value = 100112017
$ws.Range("G89").Value = "Apio"
$ws.Range("H89").Value = "Americana (o)"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = 9000
$ws.Range("N89").Value = "$/docena de matas"
$ws.Range("O89").Value = "Provincia del Elquí"
$ws.Range("P89").Value = 1500
$ws.Range("Q89").Value = 6
$ws.Range("R89").Value = "Hortaliza"
